$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: Total Trades 17 -> 18, Win Rate % 35.29 -> 33.33
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 18
$summary.Range("B9").Value = 33.33

# ---------------------------------------------------------------------------
# Strategy Status sheet: MarketMaking row - Trades 17 -> 18, Win Rate % 35.29 -> 33.33
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 18
$status.Range("G4").Value = 33.33

# ---------------------------------------------------------------------------
# New trade #18 row to append to both "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------------
function Add-TradeRow {
    param($ws, $rowIndex)

    $ws.Cells.Item($rowIndex, 1).Value = 18

    # The date column parses "YYYY-MM-DD" strings as real dates by default,
    # which would store the cell as a date serial instead of text (as the
    # rest of the column uses). Force text formatting, assign, then restore
    # the default "Normal" style so the cell matches its siblings exactly.
    $dateCell = $ws.Cells.Item($rowIndex, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Cells.Item($rowIndex, 3).Value = "12:28:57"
    $ws.Cells.Item($rowIndex, 4).Value = "MarketMaking"
    $ws.Cells.Item($rowIndex, 5).Value = "DOWN"
    $ws.Cells.Item($rowIndex, 6).Value = 0.97
    $ws.Cells.Item($rowIndex, 7).Value = 0.97
    $ws.Cells.Item($rowIndex, 8).Value = "CLOSED"
    $ws.Cells.Item($rowIndex, 9).Value = 0
    $ws.Cells.Item($rowIndex, 10).Value = 0
    $ws.Cells.Item($rowIndex, 11).Value = 99.98999999999999
    $ws.Cells.Item($rowIndex, 12).Value = 0
    $ws.Cells.Item($rowIndex, 13).Value = 0
    $ws.Cells.Item($rowIndex, 14).Value = 0.6
    $ws.Cells.Item($rowIndex, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($rowIndex, 16).Value = "early_exit"
    $ws.Cells.Item($rowIndex, 17).Value = 0.11
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $allTrades 19

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $marketMaking 19
